# Bump the published term to version 1.1.0, dated to match the new release.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $label = $ws.Cells.Item($r, 1).Text

    if ($label -eq "Version") {
        $ws.Cells.Item($r, 2).Value = "1.1.0"
    }
    elseif ($label -eq "Date") {
        $ws.Cells.Item($r, 2).Value = "2023-07-10T23:08:03+02:00"
    }
}
